$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-12 Tuesday" "2023-09-13 Wednesday"

Replace-Text "54÷7=" "69÷4="
Replace-Text "74÷7=" "29÷4="
Replace-Text "22÷2=" "52÷6="
Replace-Text "79÷3=" "46÷7="
Replace-Text "83÷7=" "64÷4="

Replace-Text "65÷2=" "64÷6="
Replace-Text "91÷5=" "90÷2="
Replace-Text "96÷7=" "31÷7="
Replace-Text "85÷3=" "80÷4="
Replace-Text "67÷9=" "79÷5="

Replace-Text "25÷2=" "64÷9="
Replace-Text "56÷2=" "60÷3="
Replace-Text "25÷4=" "29÷4="
Replace-Text "71÷4=" "61÷6="
Replace-Text "58÷7=" "63÷6="

Replace-Text "54÷3=" "78÷6="
Replace-Text "41÷3=" "41÷9="
Replace-Text "80÷5=" "77÷3="
Replace-Text "62÷6=" "47÷2="
Replace-Text "97÷5=" "89÷3="

Replace-Text "69÷9=" "92÷7="
Replace-Text "72÷4=" "29÷3="
Replace-Text "76÷6=" "12÷6="
Replace-Text "89÷9=" "95÷5="
Replace-Text "14÷8=" "91÷9="
